$wb = $excel.ActiveWorkbook

# --- Update the text note on "Hoja1" (A1) with the new exchange rate figures ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Cells.Item(1, 1)
$oldText = $cellA1.Value2
$oldSnippet = "✅ 1000 Bs = 7.1 = 28423.3 pesos`n✅ 28423.3 pesos = 7.08 = 977.9 Bs"
$newSnippet = "✅ 1000 Bs = 7.04 = 28142.15 pesos`n✅ 28142.15 pesos = 6.99 = 955.09 Bs"
$cellA1.Value2 = $oldText.Replace($oldSnippet, $newSnippet)

# --- Update the numeric rate cells on "tasas" sheet ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 142.1
$wsTasas.Range("O10").Value = 3999
$wsTasas.Range("N12").Value = 4025
$wsTasas.Range("O12").Value = 136.6
